$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.423
$ws.Range("B3").Value = 6.221
$ws.Range("E3").Value = 16.44
$ws.Range("B4").Value = 6.532999999999999
$ws.Range("E9").Value = 17.155
$ws.Range("A11").Value = -21.464
$ws.Range("A12").Value = -21.564
$ws.Range("B14").Value = 5.944
$ws.Range("A15").Value = -21.529
$ws.Range("E15").Value = 15.956
$ws.Range("E19").Value = 16.366
$ws.Range("E20").Value = 16.321
$ws.Range("E25").Value = 17.354
$ws.Range("B26").Value = 6.803
$ws.Range("A27").Value = -21.575
$ws.Range("E27").Value = 16.701
$ws.Range("A28").Value = -21.517
$ws.Range("E28").Value = 17.15
$ws.Range("E30").Value = 16.293
$ws.Range("A31").Value = -21.313
$ws.Range("B31").Value = 6.462000000000001
$ws.Range("A32").Value = -21.336
$ws.Range("E32").Value = 16.574
$ws.Range("B35").Value = 7.32
$ws.Range("A36").Value = -20.972
$ws.Range("B37").Value = 7.32
$ws.Range("A38").Value = -20.002
$ws.Range("B39").Value = 6.890000000000001
$ws.Range("B40").Value = 8.642999999999999
$ws.Range("E44").Value = 16.195
$ws.Range("B45").Value = 6.061
$ws.Range("A46").Value = -21.264
$ws.Range("E47").Value = 16.66
$ws.Range("B52").Value = 6.013000000000001
$ws.Range("A54").Value = -21.856
$ws.Range("A55").Value = -22.214
$ws.Range("A56").Value = -21.522
$ws.Range("B57").Value = 5.332
$ws.Range("E58").Value = 16.586
$ws.Range("E62").Value = 16.311
$ws.Range("A67").Value = -21.623
$ws.Range("A69").Value = -21.656
$ws.Range("A72").Value = -21.57
$ws.Range("A73").Value = -20.628
$ws.Range("E77").Value = 16.98
$ws.Range("E78").Value = 16.564
$ws.Range("B81").Value = 6.441
$ws.Range("A83").Value = -20.15
$ws.Range("B83").Value = 7.111
$ws.Range("E84").Value = 16.962
$ws.Range("A86").Value = -21.911
$ws.Range("E89").Value = 17.339
$ws.Range("A91").Value = -21.508
$ws.Range("E91").Value = 17.39
$ws.Range("E92").Value = 17.262
$ws.Range("A93").Value = -21.49
$ws.Range("E96").Value = 16.749
$ws.Range("A99").Value = -20.938
$ws.Range("B100").Value = 5.517
$ws.Range("B102").Value = 7.153
$ws.Range("E102").Value = 16.419
